$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.7066831166592121
$ws.Range("D2").Value = 4.105095489590036
$ws.Range("E2").Value = 1.43867018000876
$ws.Range("F2").Value = 0.7066831166592127

$ws.Range("C3").Value = 0.6971472239506631
$ws.Range("D3").Value = 4.248792685760669
$ws.Range("E3").Value = 1.501580499894976
$ws.Range("F3").Value = 0.6804698690085806

$ws.Range("C4").Value = 0.4246018788223884
$ws.Range("D4").Value = 1.889906636970288
$ws.Range("E4").Value = 0.233952436166282
$ws.Range("F4").Value = 0.4246018788223885

$ws.Range("C5").Value = 0.3770829766081768
$ws.Range("D5").Value = 2.045259850466442
$ws.Range("E5").Value = 0.2547892213537854
$ws.Range("F5").Value = 0.317542833205385

$ws.Range("C6").Value = 0.4907501243415704
$ws.Range("D6").Value = 0.571774120392127
$ws.Range("E6").Value = 0.4184808130584473
$ws.Range("F6").Value = 0.4907501243415697

$ws.Range("C7").Value = 0.4689316096842102
$ws.Range("D7").Value = 0.5957137074257653
$ws.Range("E7").Value = 0.439661175291062
$ws.Range("F7").Value = 0.4378967930628836

$ws.Range("C8").Value = 0.4688519577727558
$ws.Range("D8").Value = 272.4843057051052
$ws.Range("E8").Value = 74.89868450622626
$ws.Range("F8").Value = 0.4688519577727559

$ws.Range("C9").Value = 0.4496581277214421
$ws.Range("D9").Value = 282.0857255247623
$ws.Range("E9").Value = 77.72414758689595
$ws.Range("F9").Value = 0.4280222426172862

$ws.Range("C10").Value = 0.6324113483372988
$ws.Range("D10").Value = 5.638455527227638
$ws.Range("E10").Value = 1.95033372924551
$ws.Range("F10").Value = 0.6324113483372984

$ws.Range("C11").Value = 0.6172330804647074
$ws.Range("D11").Value = 5.878397254636925
$ws.Range("E11").Value = 2.047961053186884
$ws.Range("F11").Value = 0.5946897203806566
